# Generate Report for Archive
# Update status for 13256a5a-36e3-4695-9fe0-57e2bc262ba1.md and
# 49824daf-fd92-4534-8834-eb9793041682.md from "Ready for handoff" to
# "In Translation" in the Overview sheet and the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "In Translation"
$wsZh.Range("C4").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "In Translation"
$wsDe.Range("C4").Value = "In Translation"
